$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting existing D:K to E:L
$ws.Columns("D").Insert()

# Copy formatting (incl. number format/font/style) from column E into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D and update revised values for E:K across all data rows
$ws.Range("D7").Value = 43404
$ws.Range("E7").Value = 43039
$ws.Range("F7").Value = 42674
$ws.Range("G7").Value = 42308
$ws.Range("H7").Value = 41943
$ws.Range("I7").Value = 41578
$ws.Range("J7").Value = 41213
$ws.Range("K7").Value = 40847
$ws.Range("D8").Value = 20887700
$ws.Range("E8").Value = 17806700
$ws.Range("F8").Value = 16527400
$ws.Range("G8").Value = 15097800
$ws.Range("H8").Value = 14541900
$ws.Range("I8").Value = 14012000
$ws.Range("J8").Value = 12769900
$ws.Range("K8").Value = 12180300
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "NA"
$ws.Range("H9").Value = "NA"
$ws.Range("I9").Value = "NA"
$ws.Range("J9").Value = "NA"
$ws.Range("K9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "NA"
$ws.Range("H10").Value = "NA"
$ws.Range("I10").Value = "NA"
$ws.Range("J10").Value = "NA"
$ws.Range("K10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("D15").Value = -631100
$ws.Range("E15").Value = -566300
$ws.Range("F15").Value = -509000
$ws.Range("G15").Value = -434600
$ws.Range("H15").Value = -391500
$ws.Range("I15").Value = -771000
$ws.Range("J15").Value = -666800
$ws.Range("K15").Value = -634600
$ws.Range("D17").Value = 10781400
$ws.Range("E17").Value = 8291200
$ws.Range("F17").Value = 7686200
$ws.Range("G17").Value = 6799800
$ws.Range("H17").Value = 6651700
$ws.Range("I17").Value = 6523700
$ws.Range("J17").Value = 6281900
$ws.Range("K17").Value = 6082100
$ws.Range("D18").Value = 10106400
$ws.Range("E18").Value = 9515500
$ws.Range("F18").Value = 8841200
$ws.Range("G18").Value = 8297900
$ws.Range("H18").Value = 7890100
$ws.Range("I18").Value = 7488200
$ws.Range("J18").Value = 6488000
$ws.Range("K18").Value = 6098200
$ws.Range("D20").Value = -1841200
$ws.Range("E20").Value = -1868000
$ws.Range("F20").Value = -1847100
$ws.Range("G20").Value = -1550900
$ws.Range("H20").Value = -969000
$ws.Range("I20").Value = -1276300
$ws.Range("J20").Value = -565600
$ws.Range("K20").Value = -910400
$ws.Range("D21").Value = 8897600
$ws.Range("E21").Value = 8215000
$ws.Range("F21").Value = 7504200
$ws.Range("G21").Value = 7182500
$ws.Range("H21").Value = 7313400
$ws.Range("I21").Value = 6596700
$ws.Range("J21").Value = 6255000
$ws.Range("K21").Value = 5506700
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("D23").Value = 8265200
$ws.Range("E23").Value = 7647500
$ws.Range("F23").Value = 6994100
$ws.Range("G23").Value = 6747000
$ws.Range("H23").Value = 6921200
$ws.Range("I23").Value = 6211900
$ws.Range("J23").Value = 5922400
$ws.Range("K23").Value = 5187900
$ws.Range("D24").Value = 1772700
$ws.Range("E24").Value = 1513000
$ws.Range("F24").Value = 1510700
$ws.Range("G24").Value = 1379000
$ws.Range("H24").Value = 1489900
$ws.Range("I24").Value = 1292700
$ws.Range("J24").Value = 1166900
$ws.Range("K24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("D26").Value = 6492500
$ws.Range("E26").Value = 6134500
$ws.Range("F26").Value = 5483300
$ws.Range("G26").Value = 5368000
$ws.Range("H26").Value = 5431200
$ws.Range("I26").Value = 4919200
$ws.Range("J26").Value = 4755500
$ws.Range("K26").Value = 5187900
$ws.Range("D27").Value = 6222300
$ws.Range("E27").Value = 5861400
$ws.Range("F27").Value = 5199800
$ws.Range("G27").Value = 5132800
$ws.Range("H27").Value = 5147000
$ws.Range("I27").Value = 4562800
$ws.Range("J27").Value = 4427300
$ws.Range("K27").Value = 4900500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("D32").Value = 1841200
$ws.Range("E32").Value = 1868000
$ws.Range("F32").Value = 1847100
$ws.Range("G32").Value = 1550900
$ws.Range("H32").Value = 969000
$ws.Range("I32").Value = 1276300
$ws.Range("J32").Value = 565600
$ws.Range("K32").Value = 910400
$ws.Range("D33").Value = 6222300
$ws.Range("E33").Value = 5861400
$ws.Range("F33").Value = 5199800
$ws.Range("G33").Value = 5132800
$ws.Range("H33").Value = 5147000
$ws.Range("I33").Value = 4562800
$ws.Range("J33").Value = 4427300
$ws.Range("K33").Value = 4900500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("D35").Value = 6222300
$ws.Range("E35").Value = 5861400
$ws.Range("F35").Value = 5199800
$ws.Range("G35").Value = 5132800
$ws.Range("H35").Value = 5147000
$ws.Range("I35").Value = 4562800
$ws.Range("J35").Value = 4427300
$ws.Range("K35").Value = 4900500
$ws.Range("D38").Value = 43404
$ws.Range("E38").Value = 43039
$ws.Range("F38").Value = 42674
$ws.Range("G38").Value = 42308
$ws.Range("H38").Value = 41943
$ws.Range("I38").Value = 41578
$ws.Range("J38").Value = 41213
$ws.Range("K38").Value = 40847
$ws.Range("D41").Value = 6695700
$ws.Range("E41").Value = 5823400
$ws.Range("F41").Value = 5101600
$ws.Range("G41").Value = 5004100
$ws.Range("H41").Value = 4337300
$ws.Range("I41").Value = 4055200
$ws.Range("J41").Value = 4492100
$ws.Range("K41").Value = 3298800
$ws.Range("D42").Value = 219633000
$ws.Range("E42").Value = 209122000
$ws.Range("F42").Value = 209910000
$ws.Range("G42").Value = 219526000
$ws.Range("H42").Value = 216986000
$ws.Range("I42").Value = 187184000
$ws.Range("J42").Value = 190700000
$ws.Range("K42").Value = 119505000
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("D47").Value = 3609400
$ws.Range("E47").Value = 3412900
$ws.Range("F47").Value = 3199400
$ws.Range("G47").Value = 3001400
$ws.Range("H47").Value = 2575700
$ws.Range("I47").Value = 3963700
$ws.Range("J47").Value = 3542400
$ws.Range("K47").Value = 6856500
$ws.Range("D48").Value = 1997500
$ws.Range("E48").Value = 1772000
$ws.Range("F48").Value = 1875400
$ws.Range("G48").Value = 1701300
$ws.Range("H48").Value = 1690800
$ws.Range("I48").Value = 3295400
$ws.Range("J48").Value = 1681900
$ws.Range("K48").Value = 1923600
$ws.Range("D49").Value = 13186700
$ws.Range("E49").Value = 9009400
$ws.Range("F49").Value = 9035500
$ws.Range("G49").Value = 8520500
$ws.Range("H49").Value = 8100000
$ws.Range("I49").Value = 17712200
$ws.Range("J49").Value = 9038400
$ws.Range("K49").Value = 7047000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("D52").Value = 1710200
$ws.Range("E52").Value = 1465300
$ws.Range("F52").Value = 1641000
$ws.Range("G52").Value = 1649900
$ws.Range("H52").Value = 1399100
$ws.Range("I52").Value = 1540500
$ws.Range("J52").Value = 1686400
$ws.Range("K52").Value = 1700900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("D54").Value = 743089000
$ws.Range("E54").Value = 681155000
$ws.Range("F54").Value = 667010000
$ws.Range("G54").Value = 637414000
$ws.Range("H54").Value = 599585000
$ws.Range("I54").Value = 553427000
$ws.Range("J54").Value = 497165000
$ws.Range("K54").Value = 456654000
$ws.Range("D57").Value = 6992600
$ws.Range("E57").Value = 6575100
$ws.Range("F57").Value = 5944700
$ws.Range("G57").Value = 5461800
$ws.Range("H57").Value = 5565200
$ws.Range("I57").Value = 6361500
$ws.Range("J57").Value = 9278800
$ws.Range("K57").Value = 4538700
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("D59").Value = 323700
$ws.Range("E59").Value = 303600
$ws.Range("F59").Value = 436900
$ws.Range("G59").Value = 434600
$ws.Range("H59").Value = 750900
$ws.Range("I59").Value = 617700
$ws.Range("J59").Value = 660100
$ws.Range("K59").Value = 284200
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("D61").Value = 4240500
$ws.Range("E61").Value = 4416900
$ws.Range("F61").Value = 5680600
$ws.Range("G61").Value = 4600700
$ws.Range("H61").Value = 3625000
$ws.Range("I61").Value = 4346900
$ws.Range("J61").Value = 8559200
$ws.Range("K61").Value = 6857200
$ws.Range("D62").Value = 1527900
$ws.Range("E62").Value = 1367900
$ws.Range("F62").Value = 2054000
$ws.Range("G62").Value = 1217500
$ws.Range("H62").Value = 1331400
$ws.Range("I62").Value = 1329900
$ws.Range("J62").Value = 939200
$ws.Range("K62").Value = 367200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("D66").Value = 694545000
$ws.Range("E66").Value = 636478000
$ws.Range("F66").Value = 625148000
$ws.Range("G66").Value = 598701000
$ws.Range("H66").Value = 563938000
$ws.Range("I66").Value = 520497000
$ws.Range("J66").Value = 467668000
$ws.Range("K66").Value = 433038000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("D70").Value = 3113800
$ws.Range("E70").Value = 3407700
$ws.Range("F70").Value = 2674700
$ws.Range("G70").Value = 2183500
$ws.Range("H70").Value = 2183500
$ws.Range("I70").Value = 3039400
$ws.Range("J70").Value = 6525200
$ws.Range("K70").Value = 3367900
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("D72").Value = 31121400
$ws.Range("E72").Value = 28453400
$ws.Range("F72").Value = 25975900
$ws.Range("G72").Value = 23441900
$ws.Range("H72").Value = 21422100
$ws.Range("I72").Value = 18799500
$ws.Range("J72").Value = 16379300
$ws.Range("K72").Value = 14032500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("D76").Value = 45429500
$ws.Range("E76").Value = 41269500
$ws.Range("F76").Value = 39187900
$ws.Range("G76").Value = 36529600
$ws.Range("H76").Value = 33463400
$ws.Range("I76").Value = 29891100
$ws.Range("J76").Value = 22972300
$ws.Range("K76").Value = 20247500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("D80").Value = 43404
$ws.Range("E80").Value = 43039
$ws.Range("F80").Value = 42674
$ws.Range("G80").Value = 42308
$ws.Range("H80").Value = 41943
$ws.Range("I80").Value = 41578
$ws.Range("J80").Value = 41213
$ws.Range("K80").Value = 40847
$ws.Range("D81").Value = 6222300
$ws.Range("E81").Value = 5861400
$ws.Range("F81").Value = 5199800
$ws.Range("G81").Value = 5132800
$ws.Range("H81").Value = 5147000
$ws.Range("I81").Value = 4562800
$ws.Range("J81").Value = 4427300
$ws.Range("K81").Value = 4900500
$ws.Range("D83").Value = 631100
$ws.Range("E83").Value = 566300
$ws.Range("F83").Value = 509000
$ws.Range("G83").Value = 434600
$ws.Range("H83").Value = 391500
$ws.Range("I83").Value = 384000
$ws.Range("J83").Value = 331900
$ws.Range("K83").Value = 317300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("D89").Value = 13255100
$ws.Range("E89").Value = 12342000
$ws.Range("F89").Value = 3375700
$ws.Range("G89").Value = 12224400
$ws.Range("H89").Value = 3679400
$ws.Range("I89").Value = 6663700
$ws.Range("J89").Value = 4683300
$ws.Range("K89").Value = 8706400
$ws.Range("D91").Value = -309600
$ws.Range("E91").Value = "NA"
$ws.Range("F91").Value = -259000
$ws.Range("G91").Value = -209900
$ws.Range("H91").Value = -206100
$ws.Range("I91").Value = -134000
$ws.Range("J91").Value = -646700
$ws.Range("K91").Value = -285000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("D94").Value = -10214300
$ws.Range("E94").Value = -9332400
$ws.Range("F94").Value = -2049600
$ws.Range("G94").Value = -10061000
$ws.Range("H94").Value = -436100
$ws.Range("I94").Value = -2429100
$ws.Range("J94").Value = -7009700
$ws.Range("K94").Value = -4945900
$ws.Range("D96").Value = -3104800
$ws.Range("E96").Value = -2825800
$ws.Range("F96").Value = -2677700
$ws.Range("G96").Value = -2534800
$ws.Range("H96").Value = -2429800
$ws.Range("I96").Value = -4576900
$ws.Range("J96").Value = -4038100
$ws.Range("K96").Value = -3546100
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("D100").Value = -2135900
$ws.Range("E100").Value = -2182000
$ws.Range("F100").Value = -1215300
$ws.Range("G100").Value = -1723600
$ws.Range("H100").Value = -3115300
$ws.Range("I100").Value = -4747300
$ws.Range("J100").Value = 3688300
$ws.Range("K100").Value = -3284200
$ws.Range("D101").Value = -32700
$ws.Range("E101").Value = -105700
$ws.Range("F101").Value = -13400
$ws.Range("G101").Value = 227000
$ws.Range("H101").Value = 154100
$ws.Range("I101").Value = 75900
$ws.Range("J101").Value = -65500
$ws.Range("K101").Value = -45300
$ws.Range("D102").Value = 872200
$ws.Range("E102").Value = 721900
$ws.Range("F102").Value = 97500
$ws.Range("G102").Value = 666800
$ws.Range("H102").Value = 282100
$ws.Range("I102").Value = -436900
$ws.Range("J102").Value = 1296400
$ws.Range("K102").Value = 431000
